$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.350.61'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.871.88'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7201'
$ws.Range("E5").Value = '  +2.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.09'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07909'
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3092'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.40'
$ws.Range("E10").Value = '  +1.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08257'
$ws.Range("E11").Value = '  +0.95%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.879.33'
$ws.Range("E12").Value = '  -7.38%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7241'
$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.249'
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.86'
$ws.Range("E15").Value = '  +1.35%  '

$ws.Range("D16").Value = '29.349.38'
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.852'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '243.95'
$ws.Range("E18").Value = '  +2.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007829'
$ws.Range("E19").Value = '  +0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.23'
$ws.Range("E20").Value = '  -0.31%  '

$ws.Range("D21").Value = '2.117.16'
$ws.Range("E21").Value = '  -8.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.995'
$ws.Range("E23").Value = '  +4.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1603'
$ws.Range("E25").Value = '  +12.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.54'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.988'
$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.25'
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.348'
$ws.Range("E29").Value = '  -2.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.495'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.379'
$ws.Range("E31").Value = '  +1.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.104'
$ws.Range("E32").Value = '  +1.20%  '

$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.945'
$ws.Range("E34").Value = '  +1.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.189'
$ws.Range("E35").Value = '  +0.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7259'
$ws.Range("E36").Value = '  +1.69%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01856'
$ws.Range("E38").Value = '  +0.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.701'
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("D40").Value = '1.172.45'
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9039'
$ws.Range("E41").Value = '  -2.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.131'
$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.85'
$ws.Range("E43").Value = '  +2.36%  '

$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.92'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").Value = '2.013.17'
$ws.Range("E46").Value = '  -8.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5281'
$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.785'
$ws.Range("E48").Value = '  +1.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.906'
$ws.Range("E49").Value = '  +6.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.270'
$ws.Range("E50").Value = '  +0.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4284'
$ws.Range("E51").Value = '  +0.11%  '
